$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Ajout des différentes valeurs pour l'entité user (colonnes "Longueur")
$ws.Range("C8").Value = 50
$ws.Range("C9").Value = 75

# Mise à jour de la cellule sélectionnée
$ws.Range("C10").Select()
